{"js": "// Remove the trailing \"Ver no Jupiter...\" / \"\u00a9 2020 ...\" footer block\n// (and the blank paragraph that preceded it) that followed the\n// bibliography entry ending in \"Janeiro: Editora Interci\u00eancia , 2004.\"\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate the \"Ver no Jupiter ...\" paragraph and the \"\u00a9 2020 ...\" paragraph\n// by their text content, which is unambiguous in this document.\nlet jupiterIndex = -1;\nlet copyrightIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n    const text = items[i].text;\n    if (jupiterIndex === -1 && text.indexOf(\"Ver no Jupiter\") !== -1) {\n        jupiterIndex = i;\n    }\n    if (copyrightIndex === -1 && text.indexOf(\"Powered by Jekyll\") !== -1) {\n        copyrightIndex = i;\n    }\n}\n\nif (jupiterIndex !== -1 && copyrightIndex !== -1) {\n    // The blank paragraph immediately before the \"Ver no Jupiter ...\" one\n    // (if present) is also removed, matching the diff.\n    let startIndex = jupiterIndex;\n    if (startIndex > 0 && items[startIndex - 1].text === \"\") {\n        startIndex -= 1;\n    }\n\n    // Delete from the copyright paragraph back to the start paragraph so\n    // indices of not-yet-deleted items stay valid.\n    for (let i = copyrightIndex; i >= startIndex; i--) {\n        items[i].delete();\n    }\n    await context.sync();\n}\n", "ps1": "# Remove the trailing \"Ver no Jupiter...\" / \"\u00a9 2020 ...\" footer block\n# (and the blank paragraph that preceded it) that followed the\n# bibliography entry ending in \"Janeiro: Editora Interci\u00eancia , 2004.\"\n$d = $word.ActiveDocument\n\n$jupiterIndex = -1\n$copyrightIndex = -1\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $text = $d.Paragraphs.Item($i).Range.Text\n    if ($jupiterIndex -eq -1 -and $text -like \"*Ver no Jupiter*\") {\n        $jupiterIndex = $i\n    }\n    if ($copyrightIndex -eq -1 -and $text -like \"*Powered by Jekyll*\") {\n        $copyrightIndex = $i\n    }\n}\n\nif ($jupiterIndex -ne -1 -and $copyrightIndex -ne -1) {\n    $startIndex = $jupiterIndex\n    # Include the blank paragraph immediately before \"Ver no Jupiter ...\" too.\n    if ($startIndex -gt 1) {\n        $prevText = $d.Paragraphs.Item($startIndex - 1).Range.Text.Trim()\n        if ($prevText -eq \"\") {\n            $startIndex = $startIndex - 1\n        }\n    }\n\n    # Delete from the last paragraph back to the first so indices of\n    # not-yet-deleted paragraphs stay valid.\n    for ($i = $copyrightIndex; $i -ge $startIndex; $i--) {\n        $d.Paragraphs.Item($i).Range.Delete()\n    }\n}\n"}
